$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark exercises 1-4 (rows 2-5) as done (checkbox column C)
$ws.Range("C2").Value = $true
$ws.Range("C3").Value = $true
$ws.Range("C4").Value = $true
$ws.Range("C5").Value = $true

# Move the active selection to C5 to match the diff
$ws.Range("C5").Select()
